$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, matching the style of the existing header cells (e.g. G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# New data values for the Save column (plain numbers, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
